$wb = $excel.ActiveWorkbook

# Rename the "3_months" sheet to "3 months"
$ws = $wb.Worksheets.Item("3_months")
$ws.Name = "3 months"

# Update the selection on that sheet to C864 (was F864)
$ws.Range("C864").Select()
